$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.895302
$ws.Range("H2").Value = 11.685906
$ws.Range("I2").Value = 0.1607797697193069
$ws.Range("J2").Value = 0.1607797697193069
$ws.Range("M2").Value = 1.836553
$ws.Range("N2").Value = 5.509659
$ws.Range("O2").Value = 0.08501638387171169
$ws.Range("P2").Value = 0.08501638387171168
$ws.Range("Q2").Value = 7.153928574006
$ws.Range("R2").Value = 64.38535716605399
$ws.Range("S2").Value = 0.013668914621262
$ws.Range("T2").Value = 0.013668914621262

$ws.Range("G3").Value = 3.895302
$ws.Range("H3").Value = 11.685906
$ws.Range("I3").Value = 0.1607797697193069
$ws.Range("J3").Value = 0.1607797697193069
$ws.Range("O3").Value = 0.5011890686049997
$ws.Range("P3").Value = 0.5011890686049997
$ws.Range("Q3").Value = 42.173880322682
$ws.Range("R3").Value = 379.564922904138
$ws.Range("S3").Value = 0.08058106303614575
$ws.Range("T3").Value = 0.08058106303614575

$ws.Range("G4").Value = 3.895302
$ws.Range("H4").Value = 11.685906
$ws.Range("I4").Value = 0.1607797697193069
$ws.Range("J4").Value = 0.1607797697193069
$ws.Range("M4").Value = 8.938931333333334
$ws.Range("N4").Value = 26.816794
$ws.Range("O4").Value = 0.4137945475232886
$ws.Range("P4").Value = 0.4137945475232886
$ws.Range("Q4").Value = 34.819837100596
$ws.Range("R4").Value = 313.378533905364
$ws.Range("S4").Value = 0.06652979206189914
$ws.Range("T4").Value = 0.06652979206189913

$ws.Range("G5").Value = 9.844169000000001
$ws.Range("I5").Value = 0.4063210567236994
$ws.Range("J5").Value = 0.4063210567236994
$ws.Range("M5").Value = 1.836553
$ws.Range("N5").Value = 5.509659
$ws.Range("O5").Value = 0.08501638387171169
$ws.Range("P5").Value = 0.08501638387171168
$ws.Range("Q5").Value = 18.079338109457
$ws.Range("R5").Value = 162.714042985113
$ws.Range("S5").Value = 0.03454394693358157
$ws.Range("T5").Value = 0.03454394693358157

$ws.Range("G6").Value = 9.844169000000001
$ws.Range("I6").Value = 0.4063210567236994
$ws.Range("J6").Value = 0.4063210567236994
$ws.Range("O6").Value = 0.5011890686049997
$ws.Range("P6").Value = 0.5011890686049997
$ws.Range("R6").Value = 959.2327494865111
$ws.Range("S6").Value = 0.2036436719739501
$ws.Range("T6").Value = 0.2036436719739502

$ws.Range("G7").Value = 9.844169000000001
$ws.Range("I7").Value = 0.4063210567236994
$ws.Range("J7").Value = 0.4063210567236994
$ws.Range("M7").Value = 8.938931333333334
$ws.Range("N7").Value = 26.816794
$ws.Range("O7").Value = 0.4137945475232886
$ws.Range("P7").Value = 0.4137945475232886
$ws.Range("Q7").Value = 87.99635072472869
$ws.Range("R7").Value = 791.9671565225581
$ws.Range("S7").Value = 0.1681334378161677
$ws.Range("T7").Value = 0.1681334378161677

$ws.Range("G8").Value = 10.488092
$ws.Range("H8").Value = 31.464276
$ws.Range("I8").Value = 0.4328991735569938
$ws.Range("J8").Value = 0.4328991735569938
$ws.Range("M8").Value = 1.836553
$ws.Range("N8").Value = 5.509659
$ws.Range("O8").Value = 0.08501638387171169
$ws.Range("P8").Value = 0.08501638387171168
$ws.Range("Q8").Value = 19.261936826876
$ws.Range("R8").Value = 173.357431441884
$ws.Range("S8").Value = 0.03680352231686812
$ws.Range("T8").Value = 0.03680352231686812

$ws.Range("G9").Value = 10.488092
$ws.Range("H9").Value = 31.464276
$ws.Range("I9").Value = 0.4328991735569938
$ws.Range("J9").Value = 0.4328991735569938
$ws.Range("O9").Value = 0.5011890686049997
$ws.Range("P9").Value = 0.5011890686049997
$ws.Range("Q9").Value = 113.5530792789053
$ws.Range("R9").Value = 1021.977713510148
$ws.Range("S9").Value = 0.2169643335949038
$ws.Range("T9").Value = 0.2169643335949038

$ws.Range("G10").Value = 10.488092
$ws.Range("H10").Value = 31.464276
$ws.Range("I10").Value = 0.4328991735569938
$ws.Range("J10").Value = 0.4328991735569938
$ws.Range("M10").Value = 8.938931333333334
$ws.Range("N10").Value = 26.816794
$ws.Range("O10").Value = 0.4137945475232886
$ws.Range("P10").Value = 0.4137945475232886
$ws.Range("Q10").Value = 93.75233420568267
$ws.Range("R10").Value = 843.7710078511441
$ws.Range("S10").Value = 0.1791313176452218
$ws.Range("T10").Value = 0.1791313176452218

